$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Set the artist/instrument table (J22:K25) in sorted order by Artist (column J),
# also fixing the "rythm guitar" -> "rhythm guitar" typo on John's row.
$ws1.Range("J22").Value = "George"
$ws1.Range("K22").Value = "lead guitar"
$ws1.Range("J23").Value = "John☺"
$ws1.Range("K23").Value = "rhythm guitar"
$ws1.Range("J24").Value = "Paul"
$ws1.Range("K24").Value = "bass"
$ws1.Range("J25").Value = "Ringo"
$ws1.Range("K25").Value = "drums"

# Record the sort operation (data is already in its sorted order) so the
# sheet keeps a sortState remembering the last Sort applied to this range.
$sortRange = $ws1.Range("J22:K25")
$ws1.Sort.SortFields.Clear()
$ws1.Sort.SortFields.Add($ws1.Range("J22"))
$ws1.Sort.SetRange($sortRange)
$ws1.Sort.Header = 0
$ws1.Sort.Apply()

# Add new rows 20 and 21
$ws1.Range("J20").Value = 10
$ws1.Range("K20").Value = 20
$ws1.Range("L20").Value = 30

$ws1.Range("J21").Value = "Artist"
$ws1.Range("K21").Value = "Instrument"
$ws1.Range("L21").Value = "Pay"

# New VLOOKUP formulas (row 23-24)
$ws1.Range("N23").Formula = '=VLOOKUP("Pa??",J22:L25,3, FALSE)'
$ws1.Range("O23").Formula = '=VLOOKUP("Ringo",J22:L25,3,FALSE)'
$ws1.Range("P23").Formula = '=VLOOKUP("Pual",J22:L25,3,FALSE)'

$ws1.Range("N24").Formula = '=VLOOKUP("Pa",J22:L25,3)'
$ws1.Range("O24").Formula = '=VLOOKUP("Ringo",J22:L25,3)'
$ws1.Range("P24").Formula = '=VLOOKUP("Pual",J22:L25,3)'

# New HLOOKUP formulas (row 26-28)
$ws1.Range("N26").Formula = '=HLOOKUP(22,J20:L25,4)'
$ws1.Range("O26").Formula = '=HLOOKUP(30,J20:L25,4,FALSE)'
$ws1.Range("P26").Formula = '=HLOOKUP(31,J20:L25,4,FALSE)'

$ws1.Range("N27").Formula = '=HLOOKUP("Pa?",J21:L25,3, FALSE)'
$ws1.Range("O27").Formula = '=HLOOKUP("Artist",J21:L25,3,FALSE)'
$ws1.Range("P27").Formula = '=HLOOKUP("Pya",J21:L25,3,FALSE)'

$ws1.Range("N28").Formula = '=HLOOKUP("Pa",J21:L25,3)'
$ws1.Range("O28").Formula = '=HLOOKUP("Instrument",J21:L25,3)'
$ws1.Range("P28").Formula = '=HLOOKUP("Pya",J21:L25,3)'

# Selection changes (Sheet2 selected first so Sheet1 ends up the active tab)
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("C11").Select()

$ws1.Range("O27").Select()
